# updated fastq files and related metadata
#
# 1. s2cDNASampleNumber for row 37 (F37) was re-numbered from 36 -> 37.
# 2. Four new sample rows (38-41) were appended with fresh prep dates
#    ("08.14.18", "10.18.18", "11.02.18", "10.18.18"), H.BROWN as the
#    preparer and the E7420L protocol - mirroring the existing rows'
#    column layout (A-G).
# 3. The legacy bold/centered header style and the thin-bordered "protocol"
#    column style are no longer used anywhere in the sheet, so every cell
#    is reset back to the default (unstyled) look, and the oversized
#    header row reverts to the sheet's normal row height.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 37: s2cDNASampleNumber corrected from 36 to 37 ---
$ws.Range("F37").Value2 = 37

# --- Append the new sample rows (38-41) ---
$newRows = @(
    @{ Row = 38; Date = "08.14.18"; Num = 38 },
    @{ Row = 39; Date = "10.18.18"; Num = 39 },
    @{ Row = 40; Date = "11.02.18"; Num = 40 },
    @{ Row = 41; Date = "10.18.18"; Num = 41 }
)

foreach ($r in $newRows) {
    $row  = $r.Row
    $date = $r.Date
    $num  = $r.Num

    # Format the date columns as text first so "08.14.18"-style values
    # aren't auto-coerced into date serials - they must stay plain text,
    # same as every other date cell already on the sheet.
    $ws.Range("A" + $row).NumberFormat = "@"
    $ws.Range("D" + $row).NumberFormat = "@"

    $ws.Range("A" + $row).Value2 = $date
    $ws.Range("B" + $row).Value2 = "H.BROWN"
    $ws.Range("C" + $row).Value2 = $num
    $ws.Range("D" + $row).Value2 = $date
    $ws.Range("E" + $row).Value2 = "H.BROWN"
    $ws.Range("F" + $row).Value2 = $num
    $ws.Range("G" + $row).Value2 = "E7420L"
}

# --- Drop the now-unused header/border formatting across the whole sheet ---
$ws.Cells.ClearFormats()
$ws.Rows.Item(1).EntireRow.AutoFit()

# --- Match the saved view's full-column selection ---
$ws.Range("A1:H1048576").Select() | Out-Null
